$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("labelled_M_0061_14y8m_1_fa")
$ws.Activate()

# Update cell values in rows 110-127 (columns C and D) to reflect corrected
# segmentation boundaries.
$ws.Range("D110").Value = 78.135000000000005

$ws.Range("C111").Value = 78.135000000000005
$ws.Range("D111").Value = 7.8144999999999998

$ws.Range("C112").Value = 80.144999999999996
$ws.Range("D112").Value = 81.3

$ws.Range("C113").Value = 81.5
$ws.Range("D113").Value = 83.055000000000007

$ws.Range("C114").Value = 84.055000000000007
$ws.Range("D114").Value = 85.6

$ws.Range("C115").Value = 86.6
$ws.Range("D115").Value = 87.11

$ws.Range("C118").Value = 90.1
$ws.Range("D118").Value = 90.32

$ws.Range("D123").Value = 92.575000000000003

$ws.Range("D125").Value = 97.09

$ws.Range("C126").Value = 97.09
$ws.Range("D126").Value = 98.165000000000006

$ws.Range("C127").Value = 98.665000000000006
$ws.Range("D127").Value = 99.1

# Update the saved view state: scroll position and active selection, to
# match where the author was working in the sheet.
$excel.ActiveWindow.ScrollRow = 112
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G129").Select()
